# Add a new "Github" slide to the end of the deck.
#
# The deck currently ends with a "Kaynakça" (references) slide that uses the
# Title+Content layout and already has the exact paragraph formatting we
# want for the new slide (no autofit quirks, same theme, etc). Duplicating
# it - then rewriting the *original* slide's text while leaving the
# duplicate alone - gives us a new last slide ("Kaynakça", i.e. the
# untouched duplicate) while the original slide object (still at position
# 15) becomes the new "Github" slide. This mirrors exactly what happened in
# the authored edit: a new slide was appended after slide 14 and the old
# "Kaynakça" slide got pushed one position further down.

$p = $ppt.ActivePresentation

$refIndex = $p.Slides.Count
$refSlide = $p.Slides.Item($refIndex)

# Duplicate the reference ("Kaynakça") slide; the duplicate lands right
# after it and keeps all of its original content untouched.
$dup = $refSlide.Duplicate()

# Turn the original slide object into the new "Github" slide.
$titleShape = $refSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Github"

$bodyShape = $refSlide.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "https://github.com/MerveKrcan/weatherShiny"

# Drop the numbered-list bullet inherited from the "Kaynakça" slide and
# collapse the indent back to flush-left, matching the rest of the deck's
# plain (buNone) paragraphs.
$bodyRange.ParagraphFormat.Bullet.Visible = 0
$level = $bodyShape.TextFrame.Ruler.Levels.Item(1)
$level.LeftMargin = 0
$level.FirstMargin = 0

# Hyperlink the URL text to itself.
$bodyRange.ActionSettings(1).Hyperlink.Address = "https://github.com/MerveKrcan/weatherShiny"

Write-Output ("Slides: " + $p.Slides.Count)
